# Auto-generated script to update market-data cells in Tonberry_Profits workbook
# per scheduled runner data refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1831.125
$ws.Range("I19").Value = 1065.6666
$ws.Range("K19").Value = 1065.6666
$ws.Range("M19").Value = -890.6666
$ws.Range("H98").Value = 4205.357
$ws.Range("I98").Value = 3852.2727
$ws.Range("J98").Value = 5500
$ws.Range("K98").Value = 3852.2727
$ws.Range("L98").Value = 5500
$ws.Range("M98").Value = -2354.2727
$ws.Range("N98").Value = -8496
$ws.Range("H111").Value = 12500653
$ws.Range("I111").Value = 20000438
$ws.Range("J111").Value = 1010.3333
$ws.Range("K111").Value = 60001314
$ws.Range("L111").Value = 3030.9999
$ws.Range("M111").Value = -59998247
$ws.Range("N111").Value = -9164.999899999999
$ws.Range("H122").Value = 4205.357
$ws.Range("I122").Value = 3852.2727
$ws.Range("J122").Value = 5500
$ws.Range("K122").Value = 11556.8181
$ws.Range("L122").Value = 16500
$ws.Range("M122").Value = -9106.8181
$ws.Range("N122").Value = -21400
$ws.Range("H138").Value = 2185.6667
$ws.Range("I138").Value = 2117.9033
$ws.Range("J138").Value = 2258.1035
$ws.Range("K138").Value = 6353.7099
$ws.Range("L138").Value = 6774.310500000001
$ws.Range("M138").Value = -1213.7099
$ws.Range("N138").Value = -17054.3105

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2245.7078
$ws.Range("I32").Value = 1520.3164
$ws.Range("K32").Value = 1520.3164
$ws.Range("M32").Value = -1233.3164
$ws.Range("H45").Value = 1340.1923
$ws.Range("I45").Value = 1091.8334
$ws.Range("K45").Value = 1091.8334
$ws.Range("M45").Value = -714.8334
$ws.Range("H61").Value = 9724.75
$ws.Range("I61").Value = 4500
$ws.Range("K61").Value = 4500
$ws.Range("M61").Value = -4288
$ws.Range("H74").Value = 1147.9286
$ws.Range("I74").Value = 822.9091
$ws.Range("J74").Value = 2339.6667
$ws.Range("K74").Value = 822.9091
$ws.Range("L74").Value = 2339.6667
$ws.Range("M74").Value = 51.09090000000003
$ws.Range("N74").Value = -4087.6667
$ws.Range("H77").Value = 1147.9286
$ws.Range("I77").Value = 822.9091
$ws.Range("J77").Value = 2339.6667
$ws.Range("K77").Value = 4114.5455
$ws.Range("L77").Value = 11698.3335
$ws.Range("M77").Value = 253.4544999999998
$ws.Range("N77").Value = -20434.3335
$ws.Range("H132").Value = 1826.1794
$ws.Range("I132").Value = 1220.9
$ws.Range("K132").Value = 3662.7
$ws.Range("M132").Value = -1132.7
$ws.Range("H135").Value = 22447.334
$ws.Range("J135").Value = 22447.334
$ws.Range("L135").Value = 22447.334
$ws.Range("N135").Value = -32587.334
$ws.Range("H136").Value = 9724.75
$ws.Range("I136").Value = 4500
$ws.Range("K136").Value = 13500
$ws.Range("M136").Value = -10950
$ws.Range("H139").Value = 51538.8
$ws.Range("J139").Value = 51538.8
$ws.Range("L139").Value = 51538.8
$ws.Range("N139").Value = -61818.8

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 748.25
$ws.Range("I64").Value = 764.3333
$ws.Range("J64").Value = 700
$ws.Range("K64").Value = 764.3333
$ws.Range("L64").Value = 700
$ws.Range("M64").Value = -539.3333
$ws.Range("N64").Value = -1150
$ws.Range("H67").Value = 748.25
$ws.Range("I67").Value = 764.3333
$ws.Range("J67").Value = 700
$ws.Range("K67").Value = 764.3333
$ws.Range("L67").Value = 700
$ws.Range("M67").Value = 15.66669999999999
$ws.Range("N67").Value = -2260
$ws.Range("H86").Value = 89490.35000000001
$ws.Range("I86").Value = 3308.0908
$ws.Range("K86").Value = 3308.0908
$ws.Range("M86").Value = -2185.0908
$ws.Range("H89").Value = 89490.35000000001
$ws.Range("I89").Value = 3308.0908
$ws.Range("K89").Value = 16540.454
$ws.Range("M89").Value = -10924.454
$ws.Range("H94").Value = 1025.625
$ws.Range("I94").Value = 1050.8334
$ws.Range("K94").Value = 1050.8334
$ws.Range("M94").Value = -599.8334
$ws.Range("H99").Value = 1926.3
$ws.Range("I99").Value = 1794.3077
$ws.Range("K99").Value = 1794.3077
$ws.Range("M99").Value = -296.3077000000001
$ws.Range("H107").Value = 1760.6
$ws.Range("I107").Value = 1377.8572
$ws.Range("J107").Value = 2653.6667
$ws.Range("K107").Value = 1377.8572
$ws.Range("L107").Value = 2653.6667
$ws.Range("M107").Value = 542.1428000000001
$ws.Range("N107").Value = -6493.6667
$ws.Range("H134").Value = 7574.077
$ws.Range("I134").Value = 8651.35
$ws.Range("K134").Value = 25954.05
$ws.Range("M134").Value = -23419.05

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2672.1304
$ws.Range("I31").Value = 2466.5
$ws.Range("K31").Value = 2466.5
$ws.Range("M31").Value = -2171.5
$ws.Range("H34").Value = 2672.1304
$ws.Range("I34").Value = 2466.5
$ws.Range("K34").Value = 2466.5
$ws.Range("M34").Value = -2264.5
$ws.Range("H58").Value = 1892103.8
$ws.Range("I58").Value = 2071689.9
$ws.Range("J58").Value = 6450
$ws.Range("K58").Value = 2071689.9
$ws.Range("L58").Value = 6450
$ws.Range("M58").Value = -2071486.9
$ws.Range("N58").Value = -6856
$ws.Range("H99").Value = 2879.8
$ws.Range("I99").Value = 2449.5
$ws.Range("K99").Value = 2449.5
$ws.Range("M99").Value = -951.5
$ws.Range("H105").Value = 1704.4286
$ws.Range("I105").Value = 1486.6666
$ws.Range("J105").Value = 3011
$ws.Range("K105").Value = 1486.6666
$ws.Range("L105").Value = 3011
$ws.Range("M105").Value = 260.3334
$ws.Range("N105").Value = -6505
$ws.Range("H122").Value = 2462.9092
$ws.Range("I122").Value = 944.6429000000001
$ws.Range("K122").Value = 2833.9287
$ws.Range("M122").Value = -383.9287000000004
$ws.Range("H126").Value = 2879.8
$ws.Range("I126").Value = 2449.5
$ws.Range("K126").Value = 7348.5
$ws.Range("M126").Value = -4878.5
$ws.Range("H134").Value = 1273.919
$ws.Range("I134").Value = 1283.9714
$ws.Range("J134").Value = 1098
$ws.Range("K134").Value = 3851.9142
$ws.Range("L134").Value = 3294
$ws.Range("M134").Value = -1316.9142
$ws.Range("N134").Value = -8364
$ws.Range("H136").Value = 1892103.8
$ws.Range("I136").Value = 2071689.9
$ws.Range("J136").Value = 6450
$ws.Range("K136").Value = 6215069.699999999
$ws.Range("L136").Value = 19350
$ws.Range("M136").Value = -6212519.699999999
$ws.Range("N136").Value = -24450

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 89.066666
$ws.Range("J12").Value = 119.5
$ws.Range("L12").Value = 358.5
$ws.Range("N12").Value = -704.5
$ws.Range("H13").Value = 299.5
$ws.Range("I13").Value = 299.5
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 898.5
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -730.5
$ws.Range("N13").ClearContents()
$ws.Range("H56").Value = 8034.8184
$ws.Range("I56").Value = 8034.8184
$ws.Range("K56").Value = 8034.8184
$ws.Range("M56").Value = -7504.8184
$ws.Range("H113").Value = 9893.272000000001
$ws.Range("J113").Value = 869.44446
$ws.Range("L113").Value = 2608.33338
$ws.Range("N113").Value = -6948.33338
$ws.Range("H131").Value = 10007.173
$ws.Range("J131").Value = 10456.192
$ws.Range("L131").Value = 31368.576
$ws.Range("N131").Value = -41448.576

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 20000
$ws.Range("J54").Value = 20000
$ws.Range("L54").Value = 20000
$ws.Range("N54").Value = -20780
$ws.Range("H70").Value = 5050
$ws.Range("I70").Value = 5600
$ws.Range("K70").Value = 5600
$ws.Range("M70").Value = -5330
$ws.Range("H73").Value = 5050
$ws.Range("I73").Value = 5600
$ws.Range("K73").Value = 5600
$ws.Range("M73").Value = -4664
$ws.Range("H122").Value = 1783.2963
$ws.Range("I122").Value = 1430.7894
$ws.Range("J122").Value = 2620.5
$ws.Range("K122").Value = 4292.3682
$ws.Range("L122").Value = 7861.5
$ws.Range("M122").Value = -1842.3682
$ws.Range("N122").Value = -12761.5
$ws.Range("H132").Value = 1542180.6
$ws.Range("I132").Value = 2266136.5
$ws.Range("J132").Value = 3774.625
$ws.Range("K132").Value = 6798409.5
$ws.Range("L132").Value = 11323.875
$ws.Range("M132").Value = -6795879.5
$ws.Range("N132").Value = -16383.875

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2292.1738
$ws.Range("I7").Value = 1558.1428
$ws.Range("K7").Value = 1558.1428
$ws.Range("M7").Value = -1446.1428
$ws.Range("H126").Value = 2292.1738
$ws.Range("I126").Value = 1558.1428
$ws.Range("K126").Value = 4674.428400000001
$ws.Range("M126").Value = -2204.428400000001
$ws.Range("H132").Value = 1711.125
$ws.Range("I132").Value = 1373.4546
$ws.Range("K132").Value = 4120.3638
$ws.Range("M132").Value = -1590.3638
$ws.Range("H133").Value = 79663
$ws.Range("J133").Value = 79663
$ws.Range("L133").Value = 79663
$ws.Range("N133").Value = -84723

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 412.31818
$ws.Range("I113").Value = 310.4375
$ws.Range("K113").Value = 931.3125
$ws.Range("M113").Value = 1238.6875
$ws.Range("H126").Value = 8991.368
$ws.Range("I126").Value = 10082.333
$ws.Range("J126").Value = 7121.143
$ws.Range("K126").Value = 30246.999
$ws.Range("L126").Value = 21363.429
$ws.Range("M126").Value = -27776.999
$ws.Range("N126").Value = -26303.429
$ws.Range("H136").Value = 11339156
$ws.Range("I136").Value = 15016002
$ws.Range("J136").Value = 2216.6667
$ws.Range("K136").Value = 45048006
$ws.Range("L136").Value = 6650.000100000001
$ws.Range("M136").Value = -45045456
$ws.Range("N136").Value = -11750.0001
